$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STATUS")  # the "STATUS" tab (also $wb.ActiveSheet)

function Set-TextValue($range, [string]$text) {
    # Assigning a plain "NN%" string via .Value auto-converts to a numeric
    # percentage (Excel's normal typed-input behaviour), which would also
    # swap in a new percentage-formatted style. The source data stores these
    # as literal text, so force it through a text formula and then collapse
    # the formula down to its cached string result via a values-only paste;
    # that keeps the cell's existing style/format untouched.
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# Row 2: Alice -> Eric, with all stats reset to 0 / 0%
$ws.Range("A2").Value = "Eric"
Set-TextValue $ws.Range("B2") "0%"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
Set-TextValue $ws.Range("E2") "0%"
Set-TextValue $ws.Range("F2") "0%"

# Row 3: John row stats reset to 0 / 0% (name unchanged)
Set-TextValue $ws.Range("B3") "0%"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
Set-TextValue $ws.Range("E3") "0%"
Set-TextValue $ws.Range("F3") "0%"

$excel.CutCopyMode = $false
